$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05903766666666666
$ws.Range("H2").Value = 0.177113
$ws.Range("I2").Value = 0.01178879857463852
$ws.Range("J2").Value = 0.01178879857463852
$ws.Range("M2").Value = 0.02507166666666667
$ws.Range("N2").Value = 0.075215
$ws.Range("O2").Value = 0.009392568139045224
$ws.Range("P2").Value = 0.009392568139045224
$ws.Range("Q2").Value = 0.001480172699444445
$ws.Range("R2").Value = 0.013321554295
$ws.Range("S2").Value = 0.0001107270938897715
$ws.Range("T2").Value = 0.0001107270938897715
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05903766666666666
$ws.Range("H3").Value = 0.177113
$ws.Range("I3").Value = 0.01178879857463852
$ws.Range("J3").Value = 0.01178879857463852
$ws.Range("N3").Value = 7.038411000000001
$ws.Range("O3").Value = 0.8789304647757153
$ws.Range("P3").Value = 0.8789304647757155
$ws.Range("Q3").Value = 0.1385104541603333
$ws.Range("R3").Value = 1.246594087443
$ws.Range("S3").Value = 0.01036153421035432
$ws.Range("T3").Value = 0.01036153421035432
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05903766666666666
$ws.Range("H4").Value = 0.177113
$ws.Range("I4").Value = 0.01178879857463852
$ws.Range("J4").Value = 0.01178879857463852
$ws.Range("M4").Value = 0.2981003333333334
$ws.Range("N4").Value = 0.894301
$ws.Range("O4").Value = 0.1116769670852394
$ws.Range("P4").Value = 0.1116769670852394
$ws.Range("Q4").Value = 0.01759914811255556
$ws.Range("R4").Value = 0.158392333013
$ws.Range("S4").Value = 0.001316537270394423
$ws.Range("T4").Value = 0.001316537270394423
$ws.Range("G5").Value = 3.694492333333333
$ws.Range("I5").Value = 0.7377260723924206
$ws.Range("J5").Value = 0.7377260723924207
$ws.Range("M5").Value = 0.02507166666666667
$ws.Range("N5").Value = 0.075215
$ws.Range("O5").Value = 0.009392568139045224
$ws.Range("P5").Value = 0.009392568139045224
$ws.Range("Q5").Value = 0.0926270802838889
$ws.Range("R5").Value = 0.8336437225550001
$ws.Range("S5").Value = 0.00692914240289602
$ws.Range("T5").Value = 0.006929142402896021
$ws.Range("G6").Value = 3.694492333333333
$ws.Range("I6").Value = 0.7377260723924206
$ws.Range("J6").Value = 0.7377260723924207
$ws.Range("N6").Value = 7.038411000000001
$ws.Range("O6").Value = 0.8789304647757153
$ws.Range("P6").Value = 0.8789304647757155
$ws.Range("Q6").Value = 8.667785159449666
$ws.Range("S6").Value = 0.6484099196850333
$ws.Range("T6").Value = 0.6484099196850335
$ws.Range("G7").Value = 3.694492333333333
$ws.Range("I7").Value = 0.7377260723924206
$ws.Range("J7").Value = 0.7377260723924207
$ws.Range("M7").Value = 0.2981003333333334
$ws.Range("N7").Value = 0.894301
$ws.Range("O7").Value = 0.1116769670852394
$ws.Range("P7").Value = 0.1116769670852394
$ws.Range("Q7").Value = 1.101329396064111
$ws.Range("R7").Value = 9.911964564577
$ws.Range("S7").Value = 0.0823870103044913
$ws.Range("T7").Value = 0.08238701030449132
$ws.Range("G8").Value = 1.193695333333333
$ws.Range("H8").Value = 3.581086
$ws.Range("I8").Value = 0.2383602645342688
$ws.Range("J8").Value = 0.2383602645342688
$ws.Range("M8").Value = 0.02507166666666667
$ws.Range("N8").Value = 0.075215
$ws.Range("O8").Value = 0.009392568139045224
$ws.Range("P8").Value = 0.009392568139045224
$ws.Range("Q8").Value = 0.02992793149888889
$ws.Range("R8").Value = 0.26935138349
$ws.Range("S8").Value = 0.002238815026278964
$ws.Range("T8").Value = 0.002238815026278965
$ws.Range("G9").Value = 1.193695333333333
$ws.Range("H9").Value = 3.581086
$ws.Range("I9").Value = 0.2383602645342688
$ws.Range("J9").Value = 0.2383602645342688
$ws.Range("N9").Value = 7.038411000000001
$ws.Range("O9").Value = 0.8789304647757153
$ws.Range("P9").Value = 0.8789304647757155
$ws.Range("Q9").Value = 2.800572788260667
$ws.Range("R9").Value = 25.205155094346
$ws.Range("S9").Value = 0.2095020980911673
$ws.Range("T9").Value = 0.2095020980911674
$ws.Range("G10").Value = 1.193695333333333
$ws.Range("H10").Value = 3.581086
$ws.Range("I10").Value = 0.2383602645342688
$ws.Range("J10").Value = 0.2383602645342688
$ws.Range("M10").Value = 0.2981003333333334
$ws.Range("N10").Value = 0.894301
$ws.Range("O10").Value = 0.1116769670852394
$ws.Range("P10").Value = 0.1116769670852394
$ws.Range("Q10").Value = 0.3558409767651111
$ws.Range("R10").Value = 3.202568790886
$ws.Range("S10").Value = 0.0266193514168225
$ws.Range("T10").Value = 0.0266193514168225
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.06072066666666667
$ws.Range("H11").Value = 0.182162
$ws.Range("I11").Value = 0.01212486449867204
$ws.Range("J11").Value = 0.01212486449867204
$ws.Range("M11").Value = 0.02507166666666667
$ws.Range("N11").Value = 0.075215
$ws.Range("O11").Value = 0.009392568139045224
$ws.Range("P11").Value = 0.009392568139045224
$ws.Range("Q11").Value = 0.001522368314444445
$ws.Range("R11").Value = 0.01370131483
$ws.Range("S11").Value = 0.0001138836159804676
$ws.Range("T11").Value = 0.0001138836159804676
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.06072066666666667
$ws.Range("H12").Value = 0.182162
$ws.Range("I12").Value = 0.01212486449867204
$ws.Range("J12").Value = 0.01212486449867204
$ws.Range("N12").Value = 7.038411000000001
$ws.Range("O12").Value = 0.8789304647757153
$ws.Range("P12").Value = 0.8789304647757155
$ws.Range("Q12").Value = 0.1424590027313333
$ws.Range("R12").Value = 1.282131024582
$ws.Range("S12").Value = 0.01065691278916039
$ws.Range("T12").Value = 0.01065691278916039
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.06072066666666667
$ws.Range("H13").Value = 0.182162
$ws.Range("I13").Value = 0.01212486449867204
$ws.Range("J13").Value = 0.01212486449867204
$ws.Range("M13").Value = 0.2981003333333334
$ws.Range("N13").Value = 0.894301
$ws.Range("O13").Value = 0.1116769670852394
$ws.Range("P13").Value = 0.1116769670852394
$ws.Range("Q13").Value = 0.01810085097355556
$ws.Range("R13").Value = 0.162907658762
$ws.Range("S13").Value = 0.001354068093531186
$ws.Range("T13").Value = 0.001354068093531186
